$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update topic text for row 9 (C2 due-date shift: "Describing data pt 1" now references summary statistics)
$ws.Range("C9").Value = "[Describing data pt 1: Summary statistics](https://soc333-sum23.github.io/slides/07-describingpt1.html#/title-slide)"

# Row 10: topic becomes a link, and a new Exercise (G) value is added
$ws.Range("C10").Value = "[Describing data pt 2: Filtering data frames](https://soc333-sum23.github.io/slides/08-describingpt2.html#/title-slide)"
$ws.Range("G10").Value = "Summary statistics; filtering"

# Row 11: topic text expanded, and the HW1 "Due" entry is removed (due date moved)
$ws.Range("C11").Value = "Describing data pt 3: Creating new variables; plots"
$ws.Range("F11").ClearContents()

# Row 12: topic text changed, and "Component 2" due entry removed from here...
$ws.Range("C12").Value = "Describing data pt 4: Plots"
$ws.Range("F12").ClearContents()

# ...and re-added on row 13 (the due date for Component 2 moved from row 12 to row 13)
$ws.Range("F13").Value = "Component 2: Descriptive statistics"

# Update the view: scrolled so column B is the leftmost visible column, and the
# active selection moved to C13
$ws.Application.ActiveWindow.ScrollColumn = 2
$ws.Range("C13").Select()
